# by jiankong on 1113
# Rename the "质控组" (QC group) label to "北京组" (Beijing group) on both
# sheets, and update the active sheet/selection to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: column A holds the group name for rows 2-5.
$ws1.Range("A2:A5").Value = "北京组"

# Sheet2: column A holds the group name for row 2.
$ws2.Range("A2").Value = "北京组"

# Sheet2's selection moves to A2 (no longer the active tab). Do this
# first, then activate/select on Sheet1 last so Sheet1 ends up as the
# active tab (selecting a range activates its sheet).
$ws2.Activate()
$ws2.Range("A2").Select()

# Make Sheet1 the active sheet/tab, and set its selection to A2:A5.
$ws1.Activate()
$ws1.Range("A2:A5").Select()
